# Tripadvisor New Orleans shard 154 - update
#
# The workbook has two tabs, in this order:
#   1) "hotel_info"  -> 9-col header + 1 data row (the Slidell Motel 6 record)
#   2) "review_info" -> 25-col header only, no data rows
#
# The edit:
#   - tab 1 becomes "review_info": keeps the 25-col review header, loses the
#     hotel data row
#   - tab 2 becomes "hotel_info": keeps the hotel header/data, but gains a
#     new "State" column (value "Louisiana") inserted right after
#     "Hotel_Name" (i.e. before "City")
#
# i.e. the two tabs effectively trade places/names, and the relocated
# hotel_info sheet gains a State column.

$wb = $excel.ActiveWorkbook

$sheetHotel = $wb.Worksheets.Item(1)   # currently "hotel_info"
$sheetReview = $wb.Worksheets.Item(2)  # currently "review_info"

# ---------------------------------------------------------------------
# 1) Make room in tab 2 for the relocated hotel data: insert a blank
#    column at C (so B="Hotel_Name" stays put and the rest of the old
#    review_info header columns shift out of the way), then wipe that
#    sheet's old header entirely -- we'll rebuild it from scratch.
# ---------------------------------------------------------------------
$sheetReview.Columns.Item(3).Insert()
$sheetReview.Cells.Clear()

# ---------------------------------------------------------------------
# 2) Relocate the hotel header + data row from tab 1 into tab 2, using
#    Copy (not re-typed literals) so cells that were stored as text
#    (e.g. "85"/"17"/"89") keep their original text type instead of
#    being re-interpreted as numbers.
# ---------------------------------------------------------------------
$srcHeaderCols  = @("A","B","C","D","E","F","G","H","I")
$dstHeaderCols  = @("A","B","D","E","F","G","H","I","J")

for ($i = 0; $i -lt $srcHeaderCols.Length; $i++) {
    $src = $sheetHotel.Range($srcHeaderCols[$i] + "1")
    $dst = $sheetReview.Range($dstHeaderCols[$i] + "1")
    $src.Copy($dst)

    $src2 = $sheetHotel.Range($srcHeaderCols[$i] + "2")
    $dst2 = $sheetReview.Range($dstHeaderCols[$i] + "2")
    $src2.Copy($dst2)
}

# New "State" column, inserted between Hotel_Name and City.
$sheetReview.Range("C1").Value = "State"
$sheetReview.Range("C2").Value = "Louisiana"

# ---------------------------------------------------------------------
# 3) Now that tab 1's data has been relocated, wipe it and give it the
#    review_info header (25 columns, no data row).
# ---------------------------------------------------------------------
$sheetHotel.Cells.Clear()

$reviewHeaders = @(
    "STR",
    "reviewer_ID",
    "reviewer_name",
    "Review_ID",
    "Date_of_scraping",
    "ReviewURL",
    "Tripadvisor_gcode",
    "Tripadvisor_dcode",
    "Tripadvisor_rcode",
    "review_date",
    "review_title",
    "review_content",
    "review_rating",
    "trip_month",
    "trip_purpose",
    "value",
    "rooms",
    "Location",
    "Cleanliness",
    "Sleep Quality",
    "Service",
    "Picture(yes=1)",
    "respondent",
    "response_date",
    "response_text"
)

for ($i = 0; $i -lt $reviewHeaders.Length; $i++) {
    $sheetHotel.Cells.Item(1, $i + 1).Value = $reviewHeaders[$i]
}

# ---------------------------------------------------------------------
# 4) Swap the tab names (via unique placeholders so the final names
#    never collide with a currently-existing sheet name).
# ---------------------------------------------------------------------
$sheetHotel.Name = "__tmp_sheet_1__"
$sheetReview.Name = "__tmp_sheet_2__"
$sheetHotel.Name = "review_info"
$sheetReview.Name = "hotel_info"
